$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Restructure sheets: insert a new "2022-Q1" sheet right before "总计".
#    "总计" is deleted and re-created after the new sheet so that its
#    internal sheetId advances from 6 to 7, while "2022-Q1" takes the freed
#    id 6 (matching sheetId=6/7 in the target workbook). All of "总计"'s
#    rows are re-entered further down, so no information is lost.
# ---------------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item($wb.Worksheets.Count)
$oldTotal.Delete()

$lastQuarter = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastQuarter)
$wsNew.Name = "2022-Q1"

$wsTotal = $wb.Worksheets.Add($null, $wsNew)
$wsTotal.Name = "总计"

# ---------------------------------------------------------------------------
# 2. Grab a reference cell that already carries the bold/bordered/centered
#    "header + index column" style (style index 2 in the original workbook)
#    from an existing fund-holdings sheet, so the new sheet visually matches
#    its siblings.
# ---------------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item("2021-Q4").Range("B1")

# ---------------------------------------------------------------------------
# 3. Populate "2022-Q1" (fund holdings detail, same layout as the other
#    quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名).
# ---------------------------------------------------------------------------
$headers2022Q1 = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers2022Q1.Length; $c++) {
  $cell = $wsNew.Cells.Item(1, $c + 2)
  $styleSrc.Copy()
  $cell.PasteSpecial(-4122)
  $cell.Value = $headers2022Q1[$c]
}

$data2022Q1 = @(
  @("0", "501311", "嘉实恒生港股通新经济指数（LOF）A", "22.36", "94.24", "5.55", "1.2410", "7"),
  @("1", "006614", "嘉实恒生港股通新经济指数（LOF）C", "8.52", "94.24", "5.55", "0.4729", "7"),
  @("2", "006786", "泰康中证港股通大消费主题指数A", "0.85", "80.77", "5.66", "0.0481", "3"),
  @("3", "513900", "华安CES港股通精选100ETF", "2.15", "96.24", "2.11", "0.0454", "9"),
  @("4", "159735", "银华中证港股通消费主题交易型开放式指数证券投资基金", "0.72", "92.83", "6.02", "0.0433", "4"),
  @("5", "513230", "华夏中证港股通消费主题ETF", "0.60", "96.92", "6.15", "0.0369", "4"),
  @("6", "007151", "前海开源沪港深聚瑞混合", "0.60", "72.90", "5.79", "0.0347", "9"),
  @("7", "159976", "工银瑞信粤港澳大湾区创新100ETF", "1.23", "96.41", "2.39", "0.0294", "9"),
  @("8", "006787", "泰康中证港股通大消费主题指数C", "0.41", "80.77", "5.66", "0.0232", "3"),
  @("9", "513590", "鹏华中证港股通消费主题交易型开放式指数证券投资基金", "0.37", "91.21", "6.01", "0.0222", "4"),
  @("10", "517550", "招商中证沪港深消费龙头ETF", "0.34", "96.01", "5.19", "0.0176", "5"),
  @("11", "159984", "南方粤港澳大湾区创新100ETF", "0.74", "97.47", "2.34", "0.0173", "10"),
  @("12", "159983", "华夏粤港澳大湾区创新100ETF", "0.70", "98.42", "2.36", "0.0165", "10"),
  @("13", "162416", "华宝港股通恒生香港35指数(LOF)", "0.21", "94.50", "5.42", "0.0114", "5"),
  @("14", "009733", "创金合信港股通大消费精选股票A", "0.13", "82.28", "5.16", "0.0067", "5"),
  @("15", "005707", "富国港股通量化精选股票", "0.24", "80.43", "2.11", "0.0051", "8"),
  @("16", "013129", "汇添富中证沪港深消费龙头指数A", "0.11", "93.76", "4.57", "0.0050", "7"),
  @("17", "009734", "创金合信港股通大消费精选股票C", "0.07", "82.28", "5.16", "0.0036", "5"),
  @("18", "159979", "广发粤港澳大湾区创新100ETF", "0.11", "96.29", "2.86", "0.0031", "8"),
  @("19", "013130", "汇添富中证沪港深消费龙头指数C", "0.02", "93.76", "4.57", "0.0009", "7"),
)

for ($r = 0; $r -lt $data2022Q1.Length; $r++) {
  $rowVals = $data2022Q1[$r]
  $excelRow = $r + 2

  $aCell = $wsNew.Cells.Item($excelRow, 1)
  $styleSrc.Copy()
  $aCell.PasteSpecial(-4122)
  $aCell.Value = [double]$rowVals[0]

  $wsNew.Cells.Item($excelRow, 2).NumberFormat = "@"
  $wsNew.Cells.Item($excelRow, 2).Value = $rowVals[1]

  $wsNew.Cells.Item($excelRow, 3).Value = $rowVals[2]

  $wsNew.Cells.Item($excelRow, 4).NumberFormat = "@"
  $wsNew.Cells.Item($excelRow, 4).Value = $rowVals[3]

  $wsNew.Cells.Item($excelRow, 5).NumberFormat = "@"
  $wsNew.Cells.Item($excelRow, 5).Value = $rowVals[4]

  $wsNew.Cells.Item($excelRow, 6).NumberFormat = "@"
  $wsNew.Cells.Item($excelRow, 6).Value = $rowVals[5]

  $wsNew.Cells.Item($excelRow, 7).NumberFormat = "@"
  $wsNew.Cells.Item($excelRow, 7).Value = $rowVals[6]

  $wsNew.Cells.Item($excelRow, 8).Value = [double]$rowVals[7]
}

# ---------------------------------------------------------------------------
# 4. Populate "总计" (per-quarter summary: 日期/持有数量(只)/持有市值(亿元)),
#    including the new 2022-Q1 row at the top.
# ---------------------------------------------------------------------------
$headersTotal = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $headersTotal.Length; $c++) {
  $cell = $wsTotal.Cells.Item(1, $c + 2)
  $styleSrc.Copy()
  $cell.PasteSpecial(-4122)
  $cell.Value = $headersTotal[$c]
}

$dataTotal = @(
  @("0", "2022-Q1", "20", "2.08"),
  @("1", "2021-Q4", "23", "5.03"),
  @("2", "2021-Q3", "16", "4.88"),
  @("3", "2021-Q2", "17", "8.62"),
  @("4", "2021-Q1", "26", "14.11"),
  @("5", "2020-Q4", "8", "1.68"),
)

for ($r = 0; $r -lt $dataTotal.Length; $r++) {
  $rowVals = $dataTotal[$r]
  $excelRow = $r + 2

  $aCell = $wsTotal.Cells.Item($excelRow, 1)
  $styleSrc.Copy()
  $aCell.PasteSpecial(-4122)
  $aCell.Value = [double]$rowVals[0]

  $wsTotal.Cells.Item($excelRow, 2).Value = $rowVals[1]
  $wsTotal.Cells.Item($excelRow, 3).Value = [double]$rowVals[2]
  $wsTotal.Cells.Item($excelRow, 4).Value = [double]$rowVals[3]
}

Write-Host "done"
